$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 144135490.0088777
$ws.Range("C2").Value = 162163921.69962877
$ws.Range("D2").Value = 180192353.39038002
$ws.Range("E2").Value = 198220785.08113098
$ws.Range("F2").Value = 216249216.77188253

$ws.Range("B3").Value = 258460950.45082933
$ws.Range("C3").Value = 276489382.14158034
$ws.Range("D3").Value = 294517813.8323316
$ws.Range("E3").Value = 312546245.5230826
$ws.Range("F3").Value = 330574677.21383417

$ws.Range("B4").Value = 487351859.0557416
$ws.Range("C4").Value = 505380290.7464927
$ws.Range("D4").Value = 523408722.43724394
$ws.Range("E4").Value = 541437154.1279949
$ws.Range("F4").Value = 559465585.8187464

$ws.Range("B5").Value = 762443846.6366886
$ws.Range("C5").Value = 780472278.3274397
$ws.Range("D5").Value = 798500710.0181911
$ws.Range("E5").Value = 816529141.708942
$ws.Range("F5").Value = 834557573.3996935

